$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.5110453333333334
$ws.Range("H2").Value = 1.533136
$ws.Range("I2").Value = 0.1569529625135799
$ws.Range("J2").Value = 0.1569529625135799
$ws.Range("M2").Value = 0.02648366666666667
$ws.Range("N2").Value = 0.07945099999999999
$ws.Range("O2").Value = 0.001430039273477916
$ws.Range("P2").Value = 0.001430039273477917
$ws.Range("Q2").Value = 0.01353435425955556
$ws.Range("R2").Value = 0.121809188336
$ws.Range("S2").Value = 0.0002244489004831264
$ws.Range("T2").Value = 0.0002244489004831265
$ws.Range("G3").Value = 0.5110453333333334
$ws.Range("H3").Value = 1.533136
$ws.Range("I3").Value = 0.1569529625135799
$ws.Range("J3").Value = 0.1569529625135799
$ws.Range("O3").Value = 0.7016741634339546
$ws.Range("P3").Value = 0.7016741634339547
$ws.Range("Q3").Value = 6.640871253553778
$ws.Range("R3").Value = 59.767841281984
$ws.Range("S3").Value = 0.110129838670197
$ws.Range("T3").Value = 0.110129838670197
$ws.Range("G4").Value = 0.5110453333333334
$ws.Range("H4").Value = 1.533136
$ws.Range("I4").Value = 0.1569529625135799
$ws.Range("J4").Value = 0.1569529625135799
$ws.Range("O4").Value = 0.2968957972925674
$ws.Range("P4").Value = 0.2968957972925675
$ws.Range("Q4").Value = 2.809917862576
$ws.Range("R4").Value = 25.289260763184
$ws.Range("S4").Value = 0.04659867494289974
$ws.Range("T4").Value = 0.04659867494289975
$ws.Range("H5").Value = 5.654927
$ws.Range("I5").Value = 0.5789163814873767
$ws.Range("J5").Value = 0.5789163814873767
$ws.Range("M5").Value = 0.02648366666666667
$ws.Range("N5").Value = 0.07945099999999999
$ws.Range("O5").Value = 0.001430039273477916
$ws.Range("P5").Value = 0.001430039273477917
$ws.Range("Q5").Value = 0.04992106723077777
$ws.Range("R5").Value = 0.4492896050769999
$ws.Range("S5").Value = 0.0008278731615866724
$ws.Range("T5").Value = 0.0008278731615866726
$ws.Range("H6").Value = 5.654927
$ws.Range("I6").Value = 0.5789163814873767
$ws.Range("J6").Value = 0.5789163814873767
$ws.Range("O6").Value = 0.7016741634339546
$ws.Range("P6").Value = 0.7016741634339547
$ws.Range("S6").Value = 0.4062106676783672
$ws.Range("T6").Value = 0.4062106676783672
$ws.Range("H7").Value = 5.654927
$ws.Range("I7").Value = 0.5789163814873767
$ws.Range("J7").Value = 0.5789163814873767
$ws.Range("O7").Value = 0.2968957972925674
$ws.Range("P7").Value = 0.2968957972925675
$ws.Range("R7").Value = 93.278693801313
$ws.Range("S7").Value = 0.1718778406474228
$ws.Range("T7").Value = 0.1718778406474228
$ws.Range("G8").Value = 0.8600203333333334
$ws.Range("I8").Value = 0.2641306559990434
$ws.Range("J8").Value = 0.2641306559990435
$ws.Range("M8").Value = 0.02648366666666667
$ws.Range("N8").Value = 0.07945099999999999
$ws.Range("O8").Value = 0.001430039273477916
$ws.Range("P8").Value = 0.001430039273477917
$ws.Range("Q8").Value = 0.02277649183455556
$ws.Range("R8").Value = 0.204988426511
$ws.Range("S8").Value = 0.0003777172114081175
$ws.Range("T8").Value = 0.0003777172114081176
$ws.Range("G9").Value = 0.8600203333333334
$ws.Range("I9").Value = 0.2641306559990434
$ws.Range("J9").Value = 0.2641306559990435
$ws.Range("O9").Value = 0.7016741634339546
$ws.Range("P9").Value = 0.7016741634339547
$ws.Range("S9").Value = 0.1853336570853904
$ws.Range("T9").Value = 0.1853336570853905
$ws.Range("G10").Value = 0.8600203333333334
$ws.Range("I10").Value = 0.2641306559990434
$ws.Range("J10").Value = 0.2641306559990435
$ws.Range("O10").Value = 0.2968957972925674
$ws.Range("P10").Value = 0.2968957972925675
$ws.Range("R10").Value = 42.558413222259
$ws.Range("S10").Value = 0.07841928170224485
$ws.Range("T10").Value = 0.07841928170224488
